# Atualizacao de dados da ADD - metricas_retencao_anual
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: cohort_year 2020, period_index 5 -> num_customers 11 -> 12, retention_rate recalculated (num_customers / cohort_size)
$ws.Range("C7").Value = 12
$ws.Range("E7").Value = 0.1121495327102804

# Row 16: cohort_year 2022, period_index 3 -> num_customers 53 -> 54, retention_rate recalculated
$ws.Range("C16").Value = 54
$ws.Range("E16").Value = 0.2797927461139896

# Row 21: cohort_year 2024, period_index 1 -> num_customers 104 -> 105, retention_rate recalculated
$ws.Range("C21").Value = 105
$ws.Range("E21").Value = 0.5147058823529411

# Row 22: cohort_year 2025, period_index 0 -> num_customers 45 -> 51, cohort_size 45 -> 51 (retention_rate stays 1)
$ws.Range("C22").Value = 51
$ws.Range("D22").Value = 51
